# Updates cryptos list figures (Price / Volume(1h)) per the Thu Oct 12
# 05:24:12 UTC 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.824.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.561.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.485"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.59%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.91%  "

$ws.Range("E9").Value = "  -0.64%  "

$ws.Range("E10").Value = "  -1.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0864"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.782.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.571.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.76%  "

$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("E15").Value = "  -0.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.842.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("E19").Value = "  +1.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0680"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.07%  "

$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("E23").Value = "  -1.80%  "

$ws.Range("E24").Value = "  +1.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.69%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("E29").Value = "  -1.49%  "

$ws.Range("E30").Value = "  +0.78%  "

$ws.Range("E31").Value = "  -3.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.404.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.49%  "

$ws.Range("E34").Value = "  -0.80%  "

$ws.Range("E35").Value = "  -1.19%  "

$ws.Range("E36").Value = "  -0.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.916"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.12%  "

$ws.Range("E38").Value = "  -1.18%  "

$ws.Range("E39").Value = "  +1.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.808"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.68%  "

$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("E42").Value = "  +0.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.38%  "

$ws.Range("E45").Value = "  +0.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.696.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.14%  "

$ws.Range("E49").Value = "  +2.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0977"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0947"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.84%  "
